# Commit: "add rhyme element to mix"
# Insert a new row at row 124 on the active sheet for the "rhyme" TEI
# element (pushing the existing rows 124-161 down to 125-162), then
# select the newly written cell to match the author's cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert a fresh row before the current row 124; existing rows 124:161
# shift down to 125:162 automatically (dimension grows to K162).
$ws.Rows("124:124").Insert()

# Populate the new row: element name in column A, "text" in I, and
# "interpretation" in J - matching the other single-line element rows
# (e.g. row 151/158/160 "both" markers use the same I/J pairing).
$ws.Range("A124").Value = "rhyme"
$ws.Range("I124").Value = "text"
$ws.Range("J124").Value = "interpretation"

# Match the author's final selection/scroll position.
$ws.Range("I124").Select()
